$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Right" marking count (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" right count (B12): 39 -> 65
$ws.Range("B12").Value = 65

# Update correct/total marks display (E12): 24/84 -> 65/140
$ws.Range("E12").Value = "65/140"
